$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values stay as text (preserve formatting like "1.00", "0.999", etc.)
$textCells = @("D2", "D3", "D4", "D5", "D6", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D20", "D22", "D24", "D25", "D29", "D31", "D33", "D34", "D35", "D36", "D38", "D40", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Update Price (D) / Volume(1h) (E) columns for rows 2-47 ---
$ws.Range("D2").Value = "67.739.87"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "3.799.67"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "596.74"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "167.18"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "0.449"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "0.0000253"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "35.91"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "4.441.66"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "3.832.26"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "67.791.96"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "7.08"
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "461.57"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").Value = "0.700"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").Value = "83.26"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "12.08"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "3.945.52"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "29.56"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "9.05"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "3.740.23"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D44").Value = "48.10"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("D46").Value = "42.85"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  +0.07%  "

# --- Rows 48-51 reshuffled: EnergySwap moves to the top (row 48),
#     Monero/ONDO/Bittensor shift down one row, each with refreshed price/volume ---
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "27.22"
$ws.Range("E48").Value = "  +7.97%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "147.50"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.36"
$ws.Range("E50").Value = "  +10.21%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "395.74"
$ws.Range("E51").Value = "  +1.12%  "
